$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.501.19'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.053.02'
$ws.Range("E3").Value = '  -2.66%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''213.20'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").Value = '''614.04'
$ws.Range("E6").Value = '  -3.39%  '
$ws.Range("D7").Value = '''0.372'
$ws.Range("E7").Value = '  -5.43%  '
$ws.Range("E8").Value = '  +12.17%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '3.047.61'
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("D11").Value = '''0.675'
$ws.Range("E11").Value = '  +19.90%  '
$ws.Range("E12").Value = '  +5.48%  '
$ws.Range("D13").Value = '''0.0000244'
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").Value = '''5.35'
$ws.Range("D15").Value = '89.542.31'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '''32.68'
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").Value = '3.602.69'
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("D18").Value = '3.024.97'
$ws.Range("E18").Value = '  -3.98%  '
$ws.Range("D19").Value = '''3.37'
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").Value = '''13.51'
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("D22").Value = '''427.78'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").Value = '''8.33'
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = '''5.07'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("D26").Value = '''83.28'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '''11.67'
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '''1.22'
$ws.Range("E29").Value = '  +22.16%  '
$ws.Range("D30").Value = '''0.162'
$ws.Range("E30").Value = '  +2.98%  '
$ws.Range("D31").Value = '''8.60'
$ws.Range("E31").Value = '  +5.34%  '
$ws.Range("D32").Value = '''3.79'
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("D33").Value = '''505.22'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '''6.72'
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("D36").Value = '''22.86'
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("E38").Value = '  -9.06%  '
$ws.Range("D39").Value = '''22.29'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '''0.138'
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("D43").Value = '''1.84'
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").Value = '''0.361'
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").Value = '''143.33'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").Value = '''0.0699'
$ws.Range("E46").Value = '  +7.09%  '
$ws.Range("D47").Value = '''43.65'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = '''4.22'
$ws.Range("E48").Value = '  +7.15%  '
$ws.Range("D49").Value = '''161.93'
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("D50").Value = '''1.23'
$ws.Range("E50").Value = '  +3.13%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.706'
$ws.Range("E51").Value = '  -2.24%  '
